$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.066202212260725
$ws.Range("D2").Value = 1.043939020143241
$ws.Range("E2").Value = 1.070569660877392
$ws.Range("F2").Value = 1.079881621326585
$ws.Range("I2").Value = 1.047947696599273
$ws.Range("J2").Value = 1.071153419954699
$ws.Range("K2").Value = 1.046711474817277
$ws.Range("L2").Value = 1.073269147738497
$ws.Range("M2").Value = 1.082556511244861
$ws.Range("N2").Value = 1.072674580674049
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.067850654710999
$ws.Range("D3").Value = 1.044648660636677
$ws.Range("E3").Value = 1.072085904382859
$ws.Range("F3").Value = 1.081582025877659
$ws.Range("I3").Value = 1.048382716879906
$ws.Range("J3").Value = 1.072454459943789
$ws.Range("K3").Value = 1.047232517176412
$ws.Range("L3").Value = 1.074600047880616
$ws.Range("M3").Value = 1.084072927731572
$ws.Range("N3").Value = 1.073977468289155
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.068915050320177
$ws.Range("D4").Value = 1.045106854064922
$ws.Range("E4").Value = 1.073065061273071
$ws.Range("F4").Value = 1.082680587360866
$ws.Range("I4").Value = 1.048661903182632
$ws.Range("J4").Value = 1.073293693770081
$ws.Range("K4").Value = 1.047568004536366
$ws.Range("L4").Value = 1.075458769725056
$ws.Range("M4").Value = 1.085051964054109
$ws.Range("N4").Value = 1.074817893923742
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.069361993431504
$ws.Range("D5").Value = 1.045299243537194
$ws.Range("E5").Value = 1.073476240905493
$ws.Range("F5").Value = 1.083142023384114
$ws.Range("I5").Value = 1.048778725380592
$ws.Range("J5").Value = 1.073645888041431
$ws.Range("K5").Value = 1.047708648124322
$ws.Range("L5").Value = 1.075819196724272
$ws.Range("M5").Value = 1.085463037671294
$ws.Range("N5").Value = 1.075170588351355
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.06943700646988
$ws.Range("D6").Value = 1.045331532859325
$ws.Range("E6").Value = 1.073545253207844
$ws.Range("F6").Value = 1.083219477433251
$ws.Range("I6").Value = 1.048798308312943
$ws.Range("J6").Value = 1.073704986935655
$ws.Range("K6").Value = 1.047732239695007
$ws.Range("L6").Value = 1.075879680276354
$ws.Range("M6").Value = 1.085532028938325
$ws.Range("N6").Value = 1.075229771172787
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.068921024457847
$ws.Range("D7").Value = 1.045109425702672
$ws.Range("E7").Value = 1.073070557261078
$ws.Range("F7").Value = 1.082686754644904
$ws.Range("I7").Value = 1.048663466314356
$ws.Range("J7").Value = 1.073298402229343
$ws.Range("K7").Value = 1.047569885372469
$ws.Range("L7").Value = 1.075463588033252
$ws.Range("M7").Value = 1.085057458842492
$ws.Range("N7").Value = 1.074822609069556
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.066759785353592
$ws.Range("D8").Value = 1.044179053147745
$ws.Range("E8").Value = 1.071082493034013
$ws.Range("F8").Value = 1.080456642598851
$ws.Range("I8").Value = 1.048095192043594
$ws.Range("J8").Value = 1.071593661718508
$ws.Range("K8").Value = 1.046887909176881
$ws.Range("L8").Value = 1.073719446467504
$ws.Range("M8").Value = 1.083069450201991
$ws.Range("N8").Value = 1.073115447631659
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.062933603308867
$ws.Range("D9").Value = 1.042531919907417
$ws.Range("E9").Value = 1.067563859522941
$ws.Range("F9").Value = 1.076513273700065
$ws.Range("I9").Value = 1.047076057133607
$ws.Range("J9").Value = 1.068569171130612
$ws.Range("K9").Value = 1.045673335924679
$ws.Range("L9").Value = 1.070626796847273
$ws.Range("M9").Value = 1.079549119069765
$ws.Range("N9").Value = 1.070086661920215
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.060370119249275
$ws.Range("D10").Value = 1.041428505586989
$ws.Range("E10").Value = 1.06520711821793
$ws.Range("F10").Value = 1.073874497966734
$ws.Range("I10").Value = 1.046384492124143
$ws.Range("J10").Value = 1.066538474711359
$ws.Range("K10").Value = 1.044854820215757
$ws.Range("L10").Value = 1.0685515092186
$ws.Range("M10").Value = 1.07719001794081
$ws.Range("N10").Value = 1.068053081679122
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.059256928535084
$ws.Range("D11").Value = 1.040949417534922
$ws.Range("E11").Value = 1.064183877774679
$ws.Range("F11").Value = 1.072729382657914
$ws.Range("I11").Value = 1.046082112137603
$ws.Range("J11").Value = 1.065655624207226
$ws.Range("K11").Value = 1.044498269234489
$ws.Range("L11").Value = 1.067649550078663
$ws.Range("M11").Value = 1.076165459754481
$ws.Range("N11").Value = 1.067168977426009
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.05884294731399
$ws.Range("D12").Value = 1.040771264177009
$ws.Range("E12").Value = 1.063803374540294
$ws.Range("F12").Value = 1.072303645449551
$ws.Range("I12").Value = 1.045969350942309
$ws.Range("J12").Value = 1.065327150391346
$ws.Range("K12").Value = 1.044365507274224
$ws.Range("L12").Value = 1.06731400835224
$ws.Range("M12").Value = 1.075784422192497
$ws.Range("N12").Value = 1.066840037139634
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.058931770211123
$ws.Range("D13").Value = 1.040809487719307
$ws.Range("E13").Value = 1.063885013257006
$ws.Range("F13").Value = 1.072394985411915
$ws.Range("I13").Value = 1.045993558748304
$ws.Range("J13").Value = 1.065397633929529
$ws.Range("K13").Value = 1.044393999837525
$ws.Range("L13").Value = 1.06738600663047
$ws.Range("M13").Value = 1.075866177560721
$ws.Range("N13").Value = 1.066910620772526
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.059222718849579
$ws.Range("D14").Value = 1.04093469539455
$ws.Range("E14").Value = 1.064152434023043
$ws.Range("F14").Value = 1.07269419913056
$ws.Range("I14").Value = 1.046072800338459
$ws.Range("J14").Value = 1.065628483618003
$ws.Range("K14").Value = 1.044487301698815
$ws.Range("L14").Value = 1.067621824644671
$ws.Range("M14").Value = 1.076133972752599
$ws.Range("N14").Value = 1.067141798294036
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.059401916206931
$ws.Range("D15").Value = 1.041011813550501
$ws.Range("E15").Value = 1.064317143945884
$ws.Range("F15").Value = 1.072878502437566
$ws.Range("I15").Value = 1.046121564770384
$ws.Range("J15").Value = 1.065770645331864
$ws.Range("K15").Value = 1.044544745148123
$ws.Range("L15").Value = 1.067767051465181
$ws.Range("M15").Value = 1.076298907415259
$ws.Range("N15").Value = 1.067284161893837
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.060443929743194
$ws.Range("D16").Value = 1.041460273383821
$ws.Range("E16").Value = 1.065274968219738
$ws.Range("F16").Value = 1.073950441500091
$ws.Range("I16").Value = 1.04640449806787
$ws.Range("J16").Value = 1.066596990941778
$ws.Range("K16").Value = 1.044878438208975
$ws.Range("L16").Value = 1.068611297823449
$ws.Range("M16").Value = 1.077257949128055
$ws.Range("N16").Value = 1.0681116810093
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.061096694635779
$ws.Range("D17").Value = 1.041741229563619
$ws.Range("E17").Value = 1.065875039720315
$ws.Range("F17").Value = 1.074622159531784
$ws.Range("I17").Value = 1.046581188034442
$ws.Range("J17").Value = 1.067114378628567
$ws.Range("K17").Value = 1.045087182932435
$ws.Range("L17").Value = 1.069139967552728
$ws.Range("M17").Value = 1.077858704787499
$ws.Range("N17").Value = 1.06862980344595
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.061477134986654
$ws.Range("D18").Value = 1.041904981099794
$ws.Range("E18").Value = 1.066224786194911
$ws.Range("F18").Value = 1.075013720234566
$ws.Range("I18").Value = 1.046683966083575
$ws.Range("J18").Value = 1.06741582077736
$ws.Range("K18").Value = 1.0452087350756
$ws.Range("L18").Value = 1.069448009596759
$ws.Range("M18").Value = 1.078208821795196
$ws.Range("N18").Value = 1.068931673677169
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.061606803773609
$ws.Range("D19").Value = 1.041960794954913
$ws.Range("E19").Value = 1.066343996107494
$ws.Range("F19").Value = 1.075147191836584
$ws.Range("I19").Value = 1.046718962991469
$ws.Range("J19").Value = 1.067518547199241
$ws.Range("K19").Value = 1.045250146521956
$ws.Range("L19").Value = 1.069552989684298
$ws.Range("M19").Value = 1.078328153223683
$ws.Range("N19").Value = 1.069034545982353
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.061026690931297
$ws.Range("D20").Value = 1.041711098622512
$ws.Range("E20").Value = 1.065810685233419
$ws.Range("F20").Value = 1.074550115576382
$ws.Range("I20").Value = 1.046562260092937
$ws.Range("J20").Value = 1.067058903212811
$ws.Range("K20").Value = 1.045064807839805
$ws.Range("L20").Value = 1.06908327968811
$ws.Range("M20").Value = 1.07779427985021
$ws.Range("N20").Value = 1.06857424924874
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.059137055439916
$ws.Range("D21").Value = 1.040897830365453
$ws.Range("E21").Value = 1.064073697135226
$ws.Range("F21").Value = 1.072606099034829
$ws.Range("I21").Value = 1.046049477948248
$ws.Range("J21").Value = 1.065560519240988
$ws.Range("K21").Value = 1.044459835570642
$ws.Range("L21").Value = 1.067552396395377
$ws.Range("M21").Value = 1.076055126847266
$ws.Range("N21").Value = 1.06707373739981
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.057946109354056
$ws.Range("D22").Value = 1.040385345029061
$ws.Range("E22").Value = 1.062979111971102
$ws.Range("F22").Value = 1.071381553133274
$ws.Range("I22").Value = 1.045724501760079
$ws.Range("J22").Value = 1.064615273731871
$ws.Range("K22").Value = 1.044077594758031
$ws.Range("L22").Value = 1.066586890141294
$ws.Range("M22").Value = 1.074958920779702
$ws.Range("N22").Value = 1.066127149533649
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.058577727979746
$ws.Range("D23").Value = 1.040657133365532
$ws.Range("E23").Value = 1.063559610655612
$ws.Range("F23").Value = 1.072030927405552
$ws.Range("I23").Value = 1.045897022690811
$ws.Range("J23").Value = 1.065116668788624
$ws.Range("K23").Value = 1.044280406281851
$ws.Range("L23").Value = 1.067099009452084
$ws.Range("M23").Value = 1.075540303490783
$ws.Range("N23").Value = 1.066629256628894
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.061058323538706
$ws.Range("D24").Value = 1.041724713885678
$ws.Range("E24").Value = 1.065839765075052
$ws.Range("F24").Value = 1.07458266988361
$ws.Range("I24").Value = 1.04657081368705
$ws.Range("J24").Value = 1.067083971219998
$ws.Range("K24").Value = 1.045074918814976
$ws.Range("L24").Value = 1.069108895489205
$ws.Range("M24").Value = 1.077823391613463
$ws.Range("N24").Value = 1.068599352855373
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.063924947176222
$ws.Range("D25").Value = 1.042958670852036
$ws.Range("E25").Value = 1.068475400782749
$ws.Range("F25").Value = 1.077534415722696
$ws.Range("I25").Value = 1.047341653285508
$ws.Range("J25").Value = 1.069353564814236
$ws.Range("K25").Value = 1.045988871122821
$ws.Range("L25").Value = 1.071428661336122
$ws.Range("M25").Value = 1.080461314124396
$ws.Range("N25").Value = 1.070872169532841
